$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows before row 316, shifting existing rows 316:427 down to 318:429
$ws.Rows("316:317").Insert()

# Populate the two newly inserted rows (316 and 317) with new data
$ws.Range("A316").Value = 1
$ws.Range("B316").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C316").Value = "Arica y Parinacota"
$ws.Range("D316").Value = 44988
$ws.Range("E316").Value = 15
$ws.Range("F316").Value = 100114014
$ws.Range("G316").Value = "Betarraga"
$ws.Range("H316").Value = "Sin especificar"
$ws.Range("I316").Value = "Primera"
$ws.Range("J316").Value = 800
$ws.Range("K316").Value = 450
$ws.Range("L316").Value = 500
$ws.Range("M316").Value = 478
$ws.Range("N316").Value = "`$/paquete 4 unidades"
$ws.Range("O316").Value = "Región de Arica y Parinacota"
$ws.Range("P316").Value = 120
$ws.Range("Q316").Value = 4
$ws.Range("R316").Value = "Hortaliza"

$ws.Range("A317").Value = 1
$ws.Range("B317").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C317").Value = "Arica y Parinacota"
$ws.Range("D317").Value = 44988
$ws.Range("E317").Value = 15
$ws.Range("F317").Value = 100114014
$ws.Range("G317").Value = "Betarraga"
$ws.Range("H317").Value = "Sin especificar"
$ws.Range("I317").Value = "Segunda"
$ws.Range("J317").Value = 650
$ws.Range("K317").Value = 450
$ws.Range("L317").Value = 500
$ws.Range("M317").Value = 469
$ws.Range("N317").Value = "`$/paquete 5 unidades"
$ws.Range("O317").Value = "Región de Arica y Parinacota"
$ws.Range("P317").Value = 94
$ws.Range("Q317").Value = 5
$ws.Range("R317").Value = "Hortaliza"
